# Auto-generated PowerShell COM-interop edit script
# Applies numeric corrections described in the commit:
# "changed MP time limit and corrected error in fixed recourse data"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B2").Value = -583.0449769298641
$ws.Range("C2").Value = 1444.89816742
$ws.Range("C3").Value = 1484.145456856
$ws.Range("B4").Value = -591.9627421610376
$ws.Range("C4").Value = 1620.242970949
$ws.Range("C5").Value = 1320.308659059
$ws.Range("C6").Value = 1361.621325062
$ws.Range("B7").Value = -572.4375263728615
$ws.Range("C7").Value = 1485.856262025
$ws.Range("B8").Value = -572.2368907390821
$ws.Range("C8").Value = 1598.838183145
$ws.Range("C9").Value = 1374.830899549
$ws.Range("B10").Value = -580.6447985151415
$ws.Range("C10").Value = 1561.866181485
$ws.Range("B11").Value = -576.6463861113283
$ws.Range("C11").Value = 1677.300846747

$ws = $wb.Worksheets.Item("1")
$ws.Range("D2").Value = 12.70050361660852
$ws.Range("B3").Value = -583.0449769298641
$ws.Range("C3").Value = 0.09969505887898387
$ws.Range("D3").Value = 1254.2459598676494

$ws = $wb.Worksheets.Item("2")
$ws.Range("D2").Value = 10.396647928500366
$ws.Range("C3").Value = 0.08592350828960689
$ws.Range("D3").Value = 1303.0283617024304

$ws = $wb.Worksheets.Item("3")
$ws.Range("D2").Value = 11.61184968145752
$ws.Range("B3").Value = -591.9627421610376
$ws.Range("C3").Value = 0.09408212701175883
$ws.Range("D3").Value = 1436.3522895833453

$ws = $wb.Worksheets.Item("4")
$ws.Range("D2").Value = 10.428838560385499
$ws.Range("C3").Value = 0.09982818600717423
$ws.Range("D3").Value = 1129.5328513227405

$ws = $wb.Worksheets.Item("5")
$ws.Range("D2").Value = 11.105487589389893
$ws.Range("C3").Value = 0.0991460682000717
$ws.Range("D3").Value = 1176.8380528668288

$ws = $wb.Worksheets.Item("6")
$ws.Range("D2").Value = 3.2795609781677246
$ws.Range("B3").Value = -572.4375263728615
$ws.Range("C3").Value = 0.09978954906580872
$ws.Range("D3").Value = 1304.2294425332577

$ws = $wb.Worksheets.Item("7")
$ws.Range("D2").Value = 3.548194586202759
$ws.Range("B3").Value = -572.2368907390821
$ws.Range("C3").Value = 0.09798331540783724
$ws.Range("D3").Value = 1405.8753364933095

$ws = $wb.Worksheets.Item("8")
$ws.Range("D2").Value = 11.77293210564441
$ws.Range("C3").Value = 0.09911378747459428
$ws.Range("D3").Value = 1182.358108377029

$ws = $wb.Worksheets.Item("9")
$ws.Range("D2").Value = 9.715260435527588
$ws.Range("B3").Value = -580.6447985151415
$ws.Range("C3").Value = 0.09950538111434619
$ws.Range("D3").Value = 1382.1005674252017

$ws = $wb.Worksheets.Item("10")
$ws.Range("D2").Value = 10.52024122085437
$ws.Range("B3").Value = -576.6463861113283
$ws.Range("C3").Value = 0.09970735509272843
$ws.Range("D3").Value = 1486.4667990655726

